$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.685.05'
$ws.Range("E2").Value = '  +5.88%  '
$ws.Range("D3").Value = '2.656.04'
$ws.Range("E3").Value = '  +8.98%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +6.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.82'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +11.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.608'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +9.14%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +18.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.18'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +17.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '55.63'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0863'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +10.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.30'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +18.18%  '
$ws.Range("D14").Value = '3.083.79'
$ws.Range("E14").Value = '  +9.95%  '
$ws.Range("E15").Value = '  +2.87%  '
$ws.Range("D16").Value = '2.692.97'
$ws.Range("E16").Value = '  +11.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.936'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +11.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.32'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +7.09%  '
$ws.Range("D19").Value = '47.922.34'
$ws.Range("E19").Value = '  +6.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000103'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +9.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.18'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.84'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +10.57%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '282.07'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +18.05%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.97'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +8.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.09'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +11.50%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '30.46'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +42.18%  '
$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.21'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +15.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.11'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.68'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +12.04%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.32'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +4.71%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '39.86'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.19'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +14.40%  '
$ws.Range("E34").Value = '  -2.75%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.27'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +12.79%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0853'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +12.15%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.87'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.94'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.00%  '
$ws.Range("E39").Value = '  +12.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.124'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +8.53%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.16'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +14.12%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.26'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +13.80%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.70'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +40.48%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.70'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +16.41%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0335'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +14.26%  '
$ws.Range("D46").Value = '2.207.34'
$ws.Range("E46").Value = '  +10.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '96.87'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +9.87%  '
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.96'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +17.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.87'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +10.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.97'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +9.78%  '
